$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    41.74067183465279,
    41.74067389101362,
    -761.5368658848429,
    -920.4217021917362,
    -920.4217021955468,
    -1249.655897667152,
    -1281.58302684717,
    -1320.884070184373,
    -1397.511073859473,
    -1454.229251501069,
    -1503.744031488908,
    -1560.462209061932,
    -1601.209974221164,
    -1601.209974223164,
    -1689.099191146345,
    -1689.096654099501,
    -1800.509336381934,
    -1796.898307618202,
    -1796.859243203621,
    -1796.898307023876,
    -1926.252130986813,
    -1926.253501315685,
    -1926.408474484525,
    -1926.40852186969,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031,
    -1926.408521870031
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
